$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Adjusted Absorbance Material"

$headers = @(
    "Wavelength", "Material_3_minus_12.8", "Material_4_minus_12.8", "Material_5_minus_6.4", "Material_6_minus_6.4", "Material_7_minus_3.2", "Material_8_minus_3.2", "Material_9_minus_1.6", "Material_10_minus_1.6", "Material_11_minus_0.8", "Material_12_minus_0.8", "Material_13_minus_0.4", "Material_14_minus_0.4", "Material_15_minus_0.2", "Material_16_minus_0.2", "Material_17_minus_0.1", "Material_18_minus_0.1"
)
for ($i = 0; $i -lt $headers.Count; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

$data = @(
    @(630.188, -0.06046405782794295, 0.06046405782794295, 0.1398082010019932, -0.1398082010019932, 0.5006194429712831, -0.5006194429712831, 0.0508647582208468, -0.0508647582208468, -0.1305748063755403, 0.1305748063755408, -0.01876256181232439, 0.01876256181232461, -0.0101719907440796, 0.0101719907440796, -0.06455139803952892, 0.06455139803952903),
    @(710.104, -0.1048651845307083, 0.1048651845307074, -0.1708780246606452, 0.1708780246606461, 0.2920988499606114, -0.2920988499606132, 0.06204487896107125, -0.06204487896107302, -0.08940847320504708, 0.08940847320504686, -0.02390380382766111, 0.02390380382766133, 0.008070628188352691, -0.008070628188352691, -0.05850714747574881, 0.05850714747574892),
    @(800.131, -0.4884761826312021, 0.488476182631203, 0.9638394056454906, -0.9638394056454898, -0.08916514962849131, 0.0891651496284922, 0.07489134740006431, -0.07489134740006431, -0.1227386227494984, 0.1227386227494982, -0.02279500618657604, 0.02279500618657671, 0.001708332918375688, -0.001708332918375577, -0.03433295809713421, 0.03433295809713421),
    @(905.029, 0.05176736723982867, -0.05176736723982867, -0.2895843053245022, 0.2895843053245031, 0.6166226182460166, -0.6166226182460166, -0.0328982350215421, 0.03289823502154121, -0.1025515258992273, 0.1025515258992269, -0.06450248642521128, 0.06450248642521128, 0.02030056055005569, -0.02030056055005569, -0.05315768146824495, 0.05315768146824495),
    @(940.061, 0.1848390167733922, -0.1848390167733922, -0.005347875784885936, 0.005347875784885936, 0.258800990564076, -0.2588009905640769, 0.473051256782536, -0.473051256782536, 0.1946657811145065, -0.1946657811145056, -0.08324019644029823, 0.08324019644029867, -0.008327756316091284, 0.008327756316090618, -0.09476864547581987, 0.09476864547582009),
)

for ($r = 0; $r -lt $data.Count; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Count; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $row[$c]
    }
}

$headerRange = $ws.Range("A1:Q1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

Write-Output "done"